# Updates for next seminar
#  1) Refresh the "datetimeFigureOut" date placeholder shown on the slide
#     master and every slide layout: 2020-06-29 -> 2020-07-06.
#  2) On the title slide, change who questions should be directed to:
#     "Remya" -> "Homayon" (and re-flow the surrounding run so "to " gets
#     its own run, matching how PowerPoint splits runs on a partial edit).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders (slide master + all custom layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = "2020-07-06"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholders $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------
# 2) "Questions should be asked to Remya." -> "... to Homayon."
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$contentShape = $slide.Shapes.Item("Content Placeholder 4")
$textRange = $contentShape.TextFrame.TextRange

for ($pi = 1; $pi -le $textRange.Paragraphs().Count; $pi++) {
    $para = $textRange.Paragraphs($pi)
    if ($para.Text -like "Questions should be asked to Remya*") {

        $fullText = $para.Text
        $startIdx = $fullText.IndexOf("to Remya") + 1   # 1-based

        # Split "to " into its own run (same formatting as the lead-in
        # text) by touching a run-level property.
        $toRun = $para.Characters($startIdx, 3)
        $toRun.Font.Bold = $true

        # Rename the speaker.
        $nameIdx = $startIdx + 3
        $nameRun = $para.Characters($nameIdx, 5)
        $nameRun.Text = "Homayon"

        break
    }
}
